# Progress update as of 04-Nov-2025:
#   - Column H ("PERIOD TO EXPIRE") decreases by 1 day for every data row.
#   - Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025.
# Rows 3-34 on the "Training Dashboard" sheet hold the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 34; $row++) {

    # --- Column H: decrement the numeric "days to expire" counter ---
    $daysLeft = $ws.Cells.Item($row, 8).Value2
    $ws.Cells.Item($row, 8).Value = $daysLeft - 1

    # --- Column I: bump the "last update" text date to 04-Nov-2025 ---
    # Assigning the literal text "04-Nov-2025" directly gets auto-parsed
    # into a date serial by the smart-input layer (and would also bump
    # the cell style because of the added number format). Routing the
    # text through a formula and then collapsing it back to a plain
    # value via copy / paste-special-values keeps it as literal text
    # ("04-Nov-2025") in the existing style, exactly like the source cell.
    $cell = $ws.Cells.Item($row, 9)
    $cell.Formula = '="04-Nov-2025"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
